$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for rows 2-17 (columns E through T)
$data = @{
    2  = @(3, 1, 1.654227, 4.962681, 0.4107585939979205, 0.4107585939979205, 3, 1, 52.47402833333333, 157.422085, 0.699720168977827, 0.6997201689778269, 86.80395446776498, 781.2355902098849, 0.2874160728013196, 0.2874160728013195)
    3  = @(3, 1, 1.654227, 4.962681, 0.4107585939979205, 0.4107585939979205, 3, 1, 12.376362, 37.129086, 0.1650338345468634, 0.1650338345468634, 20.473312182174, 184.259809639566, 0.06778906584055507, 0.06778906584055507)
    4  = @(3, 1, 1.654227, 4.962681, 0.4107585939979205, 0.4107585939979205, 3, 1, 6.377905999999999, 19.133718, 0.08504682422503862, 0.08504682422503862, 10.550504308662, 94.95453877795799, 0.03493371394266515, 0.03493371394266515)
    5  = @(3, 1, 1.654227, 4.962681, 0.4107585939979205, 0.4107585939979205, 3, 1, 3.764580333333333, 11.293741, 0.05019917225027106, 0.05019917225027107, 6.227470431068999, 56.047233879621, 0.02061974141338077, 0.02061974141338077)
    6  = @(3, 1, 1.288726, 3.866178, 0.3200015957958394, 0.3200015957958394, 3, 1, 52.47402833333333, 157.422085, 0.699720168977827, 0.6997201689778269, 67.62464463790332, 608.6218017411298, 0.223911570683439, 0.223911570683439)
    7  = @(3, 1, 1.288726, 3.866178, 0.3200015957958394, 0.3200015957958394, 3, 1, 12.376362, 37.129086, 0.1650338345468634, 0.1650338345468634, 15.949739494812, 143.547655453308, 0.05281109041530283, 0.05281109041530283)
    8  = @(3, 1, 1.288726, 3.866178, 0.3200015957958394, 0.3200015957958394, 3, 1, 6.377905999999999, 19.133718, 0.08504682422503862, 0.08504682422503862, 8.219373287755998, 73.97435958980398, 0.02721511946938061, 0.02721511946938061)
    9  = @(3, 1, 1.288726, 3.866178, 0.3200015957958394, 0.3200015957958394, 3, 1, 3.764580333333333, 11.293741, 0.05019917225027106, 0.05019917225027107, 4.851512554655333, 43.663612991898, 0.01606381522771696, 0.01606381522771696)
    10 = @(3, 1, 0.8858993333333333, 2.657698, 0.2199763179924491, 0.2199763179924491, 3, 1, 52.47402833333333, 157.422085, 0.699720168977827, 0.6997201689778269, 46.48670671781444, 418.3803604603299, 0.1539218663967967, 0.1539218663967967)
    11 = @(3, 1, 0.8858993333333333, 2.657698, 0.2199763179924491, 0.2199763179924491, 3, 1, 12.376362, 37.129086, 0.1650338345468634, 0.1650338345468634, 10.964210844892, 98.677897604028, 0.03630353526779407, 0.03630353526779407)
    12 = @(3, 1, 0.8858993333333333, 2.657698, 0.2199763179924491, 0.2199763179924491, 3, 1, 6.377905999999999, 19.133718, 0.08504682422503862, 0.08504682422503862, 5.650182673462666, 50.851644061164, 0.01870828724997502, 0.01870828724997502)
    13 = @(3, 1, 0.8858993333333333, 2.657698, 0.2199763179924491, 0.2199763179924491, 3, 1, 3.764580333333333, 11.293741, 0.05019917225027106, 0.05019917225027107, 3.335039207579777, 30.015352868218, 0.01104262907788335, 0.01104262907788336)
    14 = @(2, 0.6666666666666666, 0.1983963333333333, 0.595189, 0.04926349221379096, 0.04926349221379096, 3, 1, 52.47402833333333, 157.422085, 0.699720168977827, 0.6997201689778269, 10.41065481656278, 93.69589334906499, 0.03447065909627167, 0.03447065909627166)
    15 = @(2, 0.6666666666666666, 0.1983963333333333, 0.595189, 0.04926349221379096, 0.04926349221379096, 3, 1, 12.376362, 37.129086, 0.1650338345468634, 0.1650338345468634, 2.455424840806, 22.098823567254, 0.008130143023211472, 0.008130143023211472)
    16 = @(2, 0.6666666666666666, 0.1983963333333333, 0.595189, 0.04926349221379096, 0.04926349221379096, 3, 1, 6.377905999999999, 19.133718, 0.08504682422503862, 0.08504682422503862, 1.265353164744667, 11.388178482702, 0.004189703563017838, 0.004189703563017838)
    17 = @(2, 0.6666666666666666, 0.1983963333333333, 0.595189, 0.04926349221379096, 0.04926349221379096, 3, 1, 3.764580333333333, 11.293741, 0.05019917225027106, 0.05019917225027107, 0.7468789346721111, 6.721910412049, 0.002472986531289979, 0.00247298653128998)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        # Column E is index 5
        $col = 5 + $i
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
